# "comit all fix of the last bugs"
# On the "Constants" sheet: insert a new row 33 containing a fresh
# Name/Value/Description constant (ALExceptionalCase / PUF / re-used
# "skip rule" description), which pushes every row below it down by one.
# Also clears the (no-op) fill formatting that had been left on B25:B26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert a brand-new row 33 (everything from the old row 33 onward shifts
# down by one row, e.g. old A34/B34 "ReviewSheet_StartCellReturnsFailedTable"
# becomes the new A35/B35, etc.).
$ws.Rows("33:33").Insert()

$ws.Range("A33").Value = "ALExceptionalCase"
$ws.Range("B33").Value = "PUF"
$ws.Range("C33").Value = "Any return which name ends with this word, must be skipped or ignored."

# Match the row height used by the rest of the table's data rows.
$ws.Rows("33:33").RowHeight = 14.25

# B25/B26 had a leftover "apply fill" style with no actual fill color;
# reset them back to the plain default style.
$ws.Range("B25:B26").Style = "Normal"

# Leave the selection on the newly-added row, like the author did.
$ws.Activate()
[void]$ws.Range("A33").Select()
